# Applies two changes to Use_case_description_Tiago.docx:
#   1. Splits the "3. The User edits the details..." run into three runs,
#      appending the parenthetical clarification about search methods.
#   2. Drops the stray <w:lastRenderedPageBreak/> before "Use case: Search Group".
#
# Both edits are performed by locating the target paragraph, taking a fresh
# Range over its *entire* content (Start..End, built via $d.Range(...) so the
# Range is not a leftover Find/Selection object), and calling InsertXML with
# the desired OOXML for that paragraph's runs. This replaces the paragraph's
# run content in place while preserving the paragraph's own attributes
# (w14:paraId, w:rsidR, <w:pPr>, ...).

$d = $word.ActiveDocument

function Get-ParagraphByText($doc, $pattern) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text -like $pattern) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# Change 1: "3. The User edits the details that he wants: ..." paragraph
# ---------------------------------------------------------------------
$p1 = Get-ParagraphByText $d "*User edits the details*"
$r1 = $d.Range($p1.Range.Start, $p1.Range.End)

$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/><w:t>3. The User edits the details that he wants: Name, Description, Color, Icon, Hierarchical Context, group Type</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> (explicit selection, searching for a keyword, free search expression, specified keywords, cited entries) </w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r>' +
  '</w:p>' +
  '</w:body>' +
  '</w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$r1.InsertXML($xml1)

# ---------------------------------------------------------------------
# Change 2: "Use case: Search Group" paragraph - drop lastRenderedPageBreak
# ---------------------------------------------------------------------
$d = $word.ActiveDocument
$p2 = Get-ParagraphByText $d "Use case: Search Group*"
$r2 = $d.Range($p2.Range.Start, $p2.Range.End)

$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p>' +
  '<w:r w:rsidRPr="009761F4"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Use case: </w:t></w:r>' +
  '<w:r w:rsidR="009B65E5"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Search</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Group</w:t></w:r>' +
  '</w:p>' +
  '</w:body>' +
  '</w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$r2.InsertXML($xml2)

Write-Host "Done."
